$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the course in row 5 (course id 7042) as done, matching the other
# rows in column C that already contain the "V" marker.
$ws.Range("C5").Value = "V"

# Move the active selection to C6, as left by the user after editing C5.
$ws.Range("C6").Select()
